$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 10
$ws.Range("W2").Value = 62
